# LED-zero35-SHIM PnP fix: correct ws2812b-mini LED orientation.
# Rotation values in column E (rows 7-66) were exported as 90 degrees but
# should be 270 degrees.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E = "Rotation" (see table1.xml). Rows 7-66 currently read 90 and
# need to become 270 (row 67 is already 270 and is left untouched).
for ($r = 7; $r -le 66; $r++) {
    $ws.Cells.Item($r, 5).Value = 270
}

# Rows 34-66 previously had no explicit cell style; give them an explicit
# "General" number format so they pick up real formatting like rows 7-33
# already had.
$ws.Range("E34:E66").NumberFormat = "General"

# Reflect the selection left behind by the edit (user had moved on to
# select the remainder of the column after fixing the first cell).
$null = $ws.Range("E8:E66").Select()
